# "New way to average"
# Lower the Powers (W) values in rows 2-4 of the "Cube 1", "Cube 3" and
# "Cube 4" sheets from 80 to 20, and update the selections left behind on
# each sheet to reflect where the user was last working.

$wb = $excel.ActiveWorkbook

# --- Cube 1 ---
$ws1 = $wb.Worksheets.Item("Cube 1")
$ws1.Range("B2:B4").Value = 20
$ws1.Range("B2:B4").Select()

# --- Cube 3 ---
$ws3 = $wb.Worksheets.Item("Cube 3")
$ws3.Range("B2:B4").Value = 20
$ws3.Range("B2:B4").Select()

# --- Cube 4 ---
$ws4 = $wb.Worksheets.Item("Cube 4")
$ws4.Range("B2:B4").Value = 20
$ws4.Range("B2:B4").Select()

# --- Cube 5 (active sheet) ---
$ws5 = $wb.Worksheets.Item("Cube 5")
$ws5.Range("D7").Select()
